# Updates cryptos list values (Price and Volume(1h) columns) for rows 2-50
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.629.65"
$ws.Range("E2").Value = "  -6.97%  "
$ws.Range("D3").Value = "1.695.25"
$ws.Range("E3").Value = "  -5.44%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'220.73"
$ws.Range("E5").Value = "  -4.74%  "
$ws.Range("D6").Value = "'0.5145"
$ws.Range("E6").Value = "  -12.49%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.2680"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").Value = "'22.17"
$ws.Range("E9").Value = "  -4.22%  "
$ws.Range("E10").Value = "  -5.66%  "
$ws.Range("D11").Value = "'0.07360"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "1.698.71"
$ws.Range("E12").Value = "  -5.29%  "
$ws.Range("D13").Value = "'4.533"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("D14").Value = "'0.5816"
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "1.926.51"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "'0.000008659"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'65.67"
$ws.Range("E17").Value = "  -12.75%  "
$ws.Range("D18").Value = "26.676.83"
$ws.Range("E18").Value = "  -6.77%  "
$ws.Range("D19").Value = "'5.022"
$ws.Range("E19").Value = "  -7.09%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'10.97"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").Value = "'187.90"
$ws.Range("E22").Value = "  -9.81%  "
$ws.Range("D23").Value = "'6.288"
$ws.Range("E23").Value = "  -7.84%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("D26").Value = "'7.534"
$ws.Range("E26").Value = "  -7.27%  "
$ws.Range("D27").Value = "'0.1182"
$ws.Range("E27").Value = "  -5.90%  "
$ws.Range("D28").Value = "'15.85"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").Value = "'1.344"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("D30").Value = "'0.05779"
$ws.Range("E30").Value = "  -7.12%  "
$ws.Range("D31").Value = "'1.342"
$ws.Range("E31").Value = "  -5.56%  "
$ws.Range("D32").Value = "'3.538"
$ws.Range("E32").Value = "  -6.38%  "
$ws.Range("D33").Value = "'3.543"
$ws.Range("E33").Value = "  -6.82%  "
$ws.Range("D34").Value = "'1.661"
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("D35").Value = "'1.027"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "'0.6017"
$ws.Range("E36").Value = "  -5.76%  "
$ws.Range("D37").Value = "'2.363"
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").Value = "'2.687"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'0.01625"
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("D40").Value = "1.097.30"
$ws.Range("E40").Value = "  -3.86%  "
$ws.Range("D41").Value = "'0.8653"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'5.863"
$ws.Range("E42").Value = "  -8.01%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'99.81"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "1.851.64"
$ws.Range("E45").Value = "  -4.93%  "
$ws.Range("D46").Value = "'0.00000000117"
$ws.Range("E46").Value = "  +4.91%  "
$ws.Range("D47").Value = "'56.61"
$ws.Range("E47").Value = "  -5.46%  "
$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "'8.156"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").Value = "'0.05251"
$ws.Range("E50").Value = "  -3.96%  "
